# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the Overview sheet's
# status text flips from "Handed back: in sync with en-US" to
# "Ready for handoff", and the associated timestamps advance a little
# (a fresh report was produced a bit later). Because the new status text
# is shorter than the old one, the status columns that were previously
# auto-sized for the long string shrink to fit the new text.
#
# Note: "2016-08-18 06:57:29" is the shared timestamp text used by both
# Overview!G2 (Latest HO Xliff Generate Date) and de-de!H2 (Latest
# Handoff Datetime) in the original report, so both move together to the
# new "2016-08-18 06:58:18" when the report is regenerated.

$wb = $excel.ActiveWorkbook

$newStatus     = "Ready for handoff"
$oldStatus     = "Handed back: in sync with en-US"

$sharedDate    = "2016-08-18 06:58:18"   # was 2016-08-18 06:57:29 (Overview!G2 & de-de!H2)
$zhHandoffDate = "2016-08-18 06:58:12"   # was 2016-08-18 06:57:24 (zh-cn!H2)

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $sharedDate

# Status text got a lot shorter, so the two status columns (E & F) that
# were auto-fit to the old text re-shrink to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $zhHandoffDate

# Status column (C) shrinks the same way.
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $sharedDate

# Status column (C) shrinks the same way.
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
